$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Expand abbreviated "Combustion Byproduct (PAH)" Table_Class entries
$pahRows = @(5, 11, 14, 18, 22, 25, 27)
foreach ($r in $pahRows) {
    $ws.Cells.Item($r, 4).Value = "Combustion Byproduct (Polycyclic Aromatic Hydrocarbon)"
}

# Expand "PFAS" Table_Class entry
$ws.Cells.Item(20, 4).Value = "Per- and Polyfluoroalkyl Substances (PFAS)"

# Correct "Insect Repellants" -> "Insect Repellents" spelling
$ws.Cells.Item(61, 4).Value = "Insecticide/Pesticide (Insect Repellents)"

$wb.Save()
